$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.061.02"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "1.828.02"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "311.72"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "0.4699"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "0.3675"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "0.07376"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "0.8773"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "20.33"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.826.63"
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("D13").Value = "0.07313"
$ws.Range("D14").Value = "5.447"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("D15").Value = "6.537"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "92.23"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "0.000008764"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "27.086.76"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "5.296"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "10.65"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "2.078.77"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("D25").Value = "1.896"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "151.70"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "18.41"
$ws.Range("D28").Value = "2.152"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("D29").Value = "5.244"
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").Value = "116.90"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "0.08900"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "0.7580"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "1.165"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "4.526"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "2.931"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "1.100"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "0.05322"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").Value = "0.01954"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "2.984"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").Value = "7.258"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "2.391"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "0.5323"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "0.1657"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "8.528"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "0.4921"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "10.48"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "103.45"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "0.06308"
$ws.Range("E51").Value = "  +0.44%  "
